$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Sample data row (row 2): India/Karnataka/Bengaluru -> Qatar/-/Doha ----
$ws.Range("A2").Value = "Qatar"
$ws.Range("C2").ClearContents()
$ws.Range("D2").Value = "Doha"

# --- New hidden helper row 3: "United Arab Emirates" ------------------------
# Styled small + white so it's effectively invisible - used as extra
# dropdown/source data for the upload template.
$a3 = $ws.Range("A3")
$a3.Font.Size = 7
$a3.Font.Color = 16777215
$a3.Font.Name = "Courier New"
$a3.Font.Family = 3
$a3.Value = "United Arab Emirates"

# --- Mobile code columns get a "+0" number format --------------------------
# Headers (already bold) first...
$ws.Range("B1").NumberFormat = "\+0"
$ws.Range("F1").NumberFormat = "\+0"

# ...then the data cells, with the new Qatar mobile code (974).
$ws.Range("B2").NumberFormat = "\+0"
$ws.Range("B2").Value = 974
$ws.Range("F2").NumberFormat = "\+0"
$ws.Range("F2").Value = 974

$ws.Range("D6").Select()
